# Refresh the regression-style summary table on Sheet1 (category, Coef.,
# Std.Err., t, P>|t|, [0.025, 0.975], coef_pos) with the updated model
# output. Most rows only change Coef. (col B) and coef_pos (col H); a few
# rows also gain/lose their Std.Err./t/P>|t|/CI columns (C:G) as the
# underlying regression now reports - or stops reporting - stats for
# that bucket.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    ,@(2, 8, 0.1011243397437618)
    ,@(3, 2, 0.1059144125693297)
    ,@(3, 8, 0.2070387523130915)
    ,@(4, 2, 0.09557263292740131)
    ,@(4, 8, 0.1966969726711631)
    ,@(5, 2, 0.05574495268222329)
    ,@(5, 8, 0.1568692924259851)
    ,@(6, 2, 0.04068131172004069)
    ,@(6, 3, $null)
    ,@(6, 4, $null)
    ,@(6, 5, $null)
    ,@(6, 6, $null)
    ,@(6, 7, $null)
    ,@(6, 8, 0.1418056514638025)
    ,@(7, 2, 0.02561311561565319)
    ,@(7, 3, $null)
    ,@(7, 4, $null)
    ,@(7, 5, $null)
    ,@(7, 6, $null)
    ,@(7, 7, $null)
    ,@(7, 8, 0.126737455359415)
    ,@(8, 2, 0.03358538300629777)
    ,@(8, 3, $null)
    ,@(8, 4, $null)
    ,@(8, 5, $null)
    ,@(8, 6, $null)
    ,@(8, 7, $null)
    ,@(8, 8, 0.1347097227500596)
    ,@(9, 2, 0.03442816579262112)
    ,@(9, 3, 0.00141594982786167)
    ,@(9, 4, 6.537405176911589)
    ,@(9, 5, 0.007345767487641948)
    ,@(9, 6, 0.03165177794909237)
    ,@(9, 7, 0.03720455363614848)
    ,@(9, 8, 0.1355525055363829)
    ,@(10, 2, 0.0339338146056304)
    ,@(10, 3, 0.002175518972394081)
    ,@(10, 4, 6.756659780268594)
    ,@(10, 5, 0.006949519492414804)
    ,@(10, 6, 0.02966593777149067)
    ,@(10, 7, 0.03820169143977017)
    ,@(10, 8, 0.1350581543493922)
    ,@(11, 2, 0.0327383655909028)
    ,@(11, 8, 0.1338627053346646)
    ,@(12, 2, 0.05449241652816593)
    ,@(12, 8, 0.1556167562719277)
    ,@(13, 2, 0.06806478292631106)
    ,@(13, 8, 0.1691891226700729)
    ,@(14, 2, 0.07758320361859486)
    ,@(14, 8, 0.1787075433623567)
    ,@(15, 2, 0.08353282968207291)
    ,@(15, 8, 0.1846571694258347)
    ,@(16, 2, 0.08650998181587267)
    ,@(16, 8, 0.1876343215596345)
    ,@(17, 2, 0.08938267756420135)
    ,@(17, 8, 0.1905070173079632)
    ,@(18, 2, -0.1011243397437618)
    ,@(18, 3, 0.008565126237189784)
    ,@(18, 4, -18.48051040529181)
    ,@(18, 5, 0.02666114851933578)
    ,@(18, 6, -0.1179468983065625)
    ,@(18, 7, -0.08430178118096128)
    ,@(19, 2, 0.08808702858097425)
    ,@(19, 8, 0.1892113683247361)
    ,@(20, 2, 0.09420210484596651)
    ,@(20, 8, 0.1953264445897283)
    ,@(21, 2, 0.09973090790117507)
    ,@(21, 3, 0.007179971999636597)
    ,@(21, 4, 24.66851722583405)
    ,@(21, 5, 0.04275798371787459)
    ,@(21, 6, 0.0856251203739136)
    ,@(21, 7, 0.113836695428437)
    ,@(21, 8, 0.2008552476449369)
    ,@(22, 2, 0.1040090445197312)
    ,@(22, 3, 0.007280550556782495)
    ,@(22, 4, 25.54785006092779)
    ,@(22, 5, 0.03768649338384446)
    ,@(22, 6, 0.08969864624144257)
    ,@(22, 7, 0.1183194427980196)
    ,@(22, 8, 0.205133384263493)
    ,@(23, 2, 0.1072261708738029)
    ,@(23, 3, 0.007381001408381161)
    ,@(23, 4, -438475234852.7534)
    ,@(23, 5, 0.03799225941035007)
    ,@(23, 6, 0.09272396700087966)
    ,@(23, 7, 0.1217283747467261)
    ,@(23, 8, 0.2083505106175647)
    ,@(24, 2, 0.1116833987350626)
    ,@(24, 3, 0.007443584060648076)
    ,@(24, 4, 26.84287002548955)
    ,@(24, 5, 0.04165186913903106)
    ,@(24, 6, 0.097049315953487)
    ,@(24, 7, 0.1263174815166377)
    ,@(24, 8, 0.2128077384788244)
    ,@(25, 2, 0.1114848809021334)
    ,@(25, 8, 0.2126092206458952)
    ,@(26, 2, 0.114026800815633)
    ,@(26, 3, 0.007233691268618695)
    ,@(26, 4, 236832989246.0349)
    ,@(26, 5, 0.04945061947872778)
    ,@(26, 6, 0.09981357336792267)
    ,@(26, 7, 0.1282400282633437)
    ,@(26, 8, 0.2151511405593948)
    ,@(27, 2, 0.1151189497054282)
    ,@(27, 8, 0.21624328944919)
    ,@(28, 2, 0.1158395932102434)
    ,@(28, 3, 0.006644235413936671)
    ,@(28, 4, 26.21812983995269)
    ,@(28, 5, 0.07409958591981149)
    ,@(28, 6, 0.1027944907970785)
    ,@(28, 7, 0.1288846956234086)
    ,@(28, 8, 0.2169639329540052)
    ,@(29, 2, 0.03391983888776901)
    ,@(29, 3, 0.001755949212353712)
    ,@(29, 4, 7.118703821210458)
    ,@(29, 5, 0.004684000656112387)
    ,@(29, 6, 0.03046947069317098)
    ,@(29, 7, 0.03737020708236766)
    ,@(29, 8, 0.1350441786315308)
)

foreach ($u in $updates) {
    $r = $u[0]
    $c = $u[1]
    $v = $u[2]
    $ws.Cells.Item($r, $c).Value = $v
}